# Updated bulk app translation sheet format to correspond to case list
# page redesign: the "module1" sheet (case list / detail properties) gets
# a new "list_or_detail" column inserted after the case_property column,
# and the "name" property row is split into two rows - one for the list
# context and one for the detail context.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("module1")

# Insert a new column before column B (default_en) to hold list_or_detail.
$ws.Columns("B").Insert()

# Header row
$ws.Range("B1").Value = "list_or_detail"

# Insert a new row below the existing "name" row (row 2) so that the
# property appears once for "list" and once for "detail".
$ws.Rows("3").Insert()

# Row 2: name / list / Name / Nom
$ws.Range("B2").Value = "list"

# Row 3: name / detail / Name / Nom (duplicate of row 2's case_property
# and translations, but for the detail context)
$ws.Range("A3").Value = "name"
$ws.Range("B3").Value = "detail"
$ws.Range("C3").Value = "Name"
$ws.Range("D3").Value = "Nom"

# Remaining original rows (other-prop, foo, baz) shifted down to rows
# 4-6, each tagged as "detail".
$ws.Range("B4").Value = "detail"
$ws.Range("B5").Value = "detail"
$ws.Range("B6").Value = "detail"

# Column widths: A & B share the old column-A width, C takes a slightly
# wider width, D keeps the old column-C width. (Inputs are pre-adjusted
# for the host's internal pixel-rounding so the saved OOXML width lands
# on the closest representable value to the target.)
$ws.Columns("A").ColumnWidth = 27.8333333
$ws.Columns("B").ColumnWidth = 27.8333333
$ws.Columns("C").ColumnWidth = 18
$ws.Columns("D").ColumnWidth = 13

# Selection within this sheet moves to A4, and this sheet becomes the
# active / selected tab (it was previously the form sheet that was
# active).
$ws.Range("A4").Select()
$ws.Activate()
